# Auto-generated edit script: apply scheduled-runner market-price updates
# to the leve-profit tables across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1109.625
$ws.Range("I98").Value = 842.4286
$ws.Range("K98").Value = 842.4286
$ws.Range("M98").Value = 655.5714
$ws.Range("H122").Value = 1109.625
$ws.Range("I122").Value = 842.4286
$ws.Range("K122").Value = 2527.2858
$ws.Range("M122").Value = -77.28579999999965
$ws.Range("H138").Value = 3261.83
$ws.Range("I138").Value = 2039.1765
$ws.Range("J138").Value = 3839.1943
$ws.Range("K138").Value = 6117.529500000001
$ws.Range("L138").Value = 11517.5829
$ws.Range("M138").Value = -977.5295000000006
$ws.Range("N138").Value = -21797.5829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3309
$ws.Range("I63").Value = 2098.3333
$ws.Range("J63").Value = 5125
$ws.Range("K63").Value = 2098.3333
$ws.Range("L63").Value = 5125
$ws.Range("M63").Value = -1412.3333
$ws.Range("N63").Value = -6497
$ws.Range("H66").Value = 3309
$ws.Range("I66").Value = 2098.3333
$ws.Range("J66").Value = 5125
$ws.Range("K66").Value = 10491.6665
$ws.Range("L66").Value = 25625
$ws.Range("M66").Value = -7059.666499999999
$ws.Range("N66").Value = -32489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 462.23077
$ws.Range("I64").Value = 218.66667
$ws.Range("J64").Value = 535.3
$ws.Range("K64").Value = 218.66667
$ws.Range("L64").Value = 535.3
$ws.Range("M64").Value = 6.333329999999989
$ws.Range("N64").Value = -985.3
$ws.Range("H67").Value = 462.23077
$ws.Range("I67").Value = 218.66667
$ws.Range("J67").Value = 535.3
$ws.Range("K67").Value = 218.66667
$ws.Range("L67").Value = 535.3
$ws.Range("M67").Value = 561.3333299999999
$ws.Range("N67").Value = -2095.3
$ws.Range("H80").Value = 160.47058
$ws.Range("I80").Value = 36
$ws.Range("J80").Value = 177.06667
$ws.Range("K80").Value = 36
$ws.Range("L80").Value = 177.06667
$ws.Range("M80").Value = 962
$ws.Range("N80").Value = -2173.06667
$ws.Range("H83").Value = 160.47058
$ws.Range("I83").Value = 36
$ws.Range("J83").Value = 177.06667
$ws.Range("K83").Value = 180
$ws.Range("L83").Value = 885.3333499999999
$ws.Range("M83").Value = 4812
$ws.Range("N83").Value = -10869.33335
$ws.Range("H94").Value = 9319.833000000001
$ws.Range("I94").Value = 859.625
$ws.Range("J94").Value = 26240.25
$ws.Range("K94").Value = 859.625
$ws.Range("L94").Value = 26240.25
$ws.Range("M94").Value = -408.625
$ws.Range("N94").Value = -27142.25
$ws.Range("H105").Value = 1978.2142
$ws.Range("I105").Value = 1329.5
$ws.Range("K105").Value = 1329.5
$ws.Range("M105").Value = 417.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 39386.883
$ws.Range("J4").Value = 10598.5625
$ws.Range("L4").Value = 10598.5625
$ws.Range("N4").Value = -10822.5625
$ws.Range("H22").Value = 267.5
$ws.Range("I22").Value = 256.66666
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 256.66666
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 93.33334000000002
$ws.Range("N22").Value = -1000
$ws.Range("H99").Value = 1626782.2
$ws.Range("I99").Value = 2980184.2
$ws.Range("J99").Value = 2699.8
$ws.Range("K99").Value = 2980184.2
$ws.Range("L99").Value = 2699.8
$ws.Range("M99").Value = -2978686.2
$ws.Range("N99").Value = -5695.8
$ws.Range("H126").Value = 1626782.2
$ws.Range("I126").Value = 2980184.2
$ws.Range("J126").Value = 2699.8
$ws.Range("K126").Value = 8940552.600000001
$ws.Range("L126").Value = 8099.400000000001
$ws.Range("M126").Value = -8938082.600000001
$ws.Range("N126").Value = -13039.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 100
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -132
$ws.Range("N13").ClearContents()
$ws.Range("H131").Value = 854.625
$ws.Range("I131").Value = 337.625
$ws.Range("J131").Value = 1371.625
$ws.Range("K131").Value = 1012.875
$ws.Range("L131").Value = 4114.875
$ws.Range("M131").Value = 4027.125
$ws.Range("N131").Value = -14194.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 118490840
$ws.Range("I70").Value = 276472000
$ws.Range("J70").Value = 4975
$ws.Range("K70").Value = 276472000
$ws.Range("L70").Value = 4975
$ws.Range("M70").Value = -276471730
$ws.Range("N70").Value = -5515
$ws.Range("H73").Value = 118490840
$ws.Range("I73").Value = 276472000
$ws.Range("J73").Value = 4975
$ws.Range("K73").Value = 276472000
$ws.Range("L73").Value = 4975
$ws.Range("M73").Value = -276471064
$ws.Range("N73").Value = -6847
$ws.Range("H132").Value = 2423.8696
$ws.Range("I132").Value = 1683.6154
$ws.Range("J132").Value = 3386.2
$ws.Range("K132").Value = 5050.8462
$ws.Range("L132").Value = 10158.6
$ws.Range("M132").Value = -2520.8462
$ws.Range("N132").Value = -15218.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 833366.7
$ws.Range("J2").Value = 833366.7
$ws.Range("L2").Value = 833366.7
$ws.Range("N2").Value = -833590.7
$ws.Range("H16").Value = 1200
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -1030
$ws.Range("H68").Value = 26028108
$ws.Range("I68").Value = 56389668
$ws.Range("J68").Value = 3912.5715
$ws.Range("K68").Value = 56389668
$ws.Range("L68").Value = 3912.5715
$ws.Range("M68").Value = -56388919
$ws.Range("N68").Value = -5410.5715
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16622
$ws.Range("H71").Value = 26028108
$ws.Range("I71").Value = 56389668
$ws.Range("J71").Value = 3912.5715
$ws.Range("K71").Value = 281948340
$ws.Range("L71").Value = 19562.8575
$ws.Range("M71").Value = -281944596
$ws.Range("N71").Value = -27050.8575
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53112

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1400
$ws.Range("I17").Value = 1400
$ws.Range("K17").Value = 1400
$ws.Range("M17").Value = -1228
$ws.Range("H55").Value = 1001
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1001
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1001
$ws.Range("N55").Value = -1555
$ws.Range("M55").ClearContents()
$ws.Range("H62").Value = 108039.93
$ws.Range("I62").Value = 4559.0435
$ws.Range("J62").Value = 504716.66
$ws.Range("K62").Value = 4559.0435
$ws.Range("L62").Value = 504716.66
$ws.Range("M62").Value = -3935.0435
$ws.Range("N62").Value = -505964.66
$ws.Range("H65").Value = 108039.93
$ws.Range("I65").Value = 4559.0435
$ws.Range("J65").Value = 504716.66
$ws.Range("K65").Value = 22795.2175
$ws.Range("L65").Value = 2523583.3
$ws.Range("M65").Value = -19675.2175
$ws.Range("N65").Value = -2529823.3
